# Actualización automática 2025-09-11 16:40:08
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("M22").Value = 467.04

$ws1.Range("H24").Value = 984.6
$ws1.Range("I24").Value = 259.2

$ws1.Range("I51").Value = 392.4

$ws1.Range("H55").Value = "3 de 53"
$ws1.Range("I55").Value = "9 de 53"

# --- Sheet "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("F22").Value = 1250.78
$ws2.Range("F24").Value = 1243.8
$ws2.Range("F52").Value = 483.98
$ws2.Range("F53").Value = 483.98
$ws2.Range("F59").Value = 28820.73

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Widen column E slightly (23 -> 24)
# (the COM width model adds a fixed 5/6 offset when round-tripping to the
# stored OOXML column width, so back it out to land exactly on 24)
$ws3.Columns.Item(5).ColumnWidth = 23.166666666666668

$ws3.Range("D6").Value = 2911.5
$ws3.Range("E6").Value = -3.916318539740132
$ws3.Range("F6").Value = 1.001346932356483

$ws3.Range("D7").Value = 2274.3
$ws3.Range("E7").Value = -1387.588983712426
$ws3.Range("F7").Value = 2.564871709299267

$ws3.Range("D12").Value = 13728.46
$ws3.Range("E12").Value = 48135.2603947566
$ws3.Range("F12").Value = 0.221914555290205

$ws3.Range("D15").Value = 27907.92
$ws3.Range("E15").Value = 94146.91551083435
$ws3.Range("F15").Value = 0.2286506706858223
